$wb = $excel.ActiveWorkbook

# The "optimization_parameters" sheet has a leftover debug/test row (row 16:
# A16="Sheet", B16=3, C16=4) sitting between the "Sigmoid" row (15) and the
# "threshold_b" row (17). Clean it up by deleting the entire row, which
# shifts everything below it up by one.
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Activate()
$ws.Rows.Item(16).Select()
$ws.Rows.Item(16).Delete()

# Move focus to the "wt_log2_expression" sheet, which becomes the active tab.
$ws2 = $wb.Worksheets.Item("wt_log2_expression")
$ws2.Activate()
